# Add a new "Date Created (Year)*" column (G) to Sheet1, with a black-font
# style applied to the header and the three data rows, then move the
# selection onto the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header + values for column G.
$ws.Range("G1").Value = "Date Created (Year)*"
$ws.Range("G2").Value = 2000
$ws.Range("G3").Value = 2000
$ws.Range("G4").Value = 2000

# Explicit black font color on the new column forces a new font/style entry
# (fontId 1 / cellXfs index 1) distinct from the default theme-colored font.
$ws.Range("G1:G4").Font.Color = 0

# Match the author's final selection state.
$ws.Range("G1:G4").Select()
